# Saldo_guide.xlsx update
# - Every existing data row (2-310) gets its "Dt. Referencia" (column G) date
#   bumped by one day (serial 45411 -> 45412).
# - Row 67 and row 121 get new "Saldo Previsto" (D) / "Vl. Total" (H) amounts.
# - A brand-new record for LYLLE MARIA LEITE PUGLIESE is inserted as row 296,
#   pushing the old rows 296-310 down to 297-311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the reference date for every existing data row by one day.
$ws.Range("G2:G310").Value2 = 45412

# Updated balances for two existing accounts.
$ws.Cells.Item(67, 4).Value2 = 73948.72
$ws.Cells.Item(67, 8).Value2 = 73948.72

$ws.Cells.Item(121, 4).Value2 = 4711.49
$ws.Cells.Item(121, 8).Value2 = 4711.49

# Insert the new record as row 296 (shifts old 296..310 -> 297..311).
$ws.Rows.Item(296).Insert()

$ws.Cells.Item(296, 1).Value2 = 15
$ws.Cells.Item(296, 2).Value2 = 806458
$ws.Cells.Item(296, 3).Value2 = "LYLLE MARIA LEITE PUGLIESE"
$ws.Cells.Item(296, 4).Value2 = 0
$ws.Cells.Item(296, 5).Value2 = 0
$ws.Cells.Item(296, 6).Value2 = 55091512772
$ws.Cells.Item(296, 7).Value2 = 45412
$ws.Cells.Item(296, 8).Value2 = 0
